$d = $word.ActiveDocument

# --- Paragraph 1 formatting: add a 5-pt-space box border and widen the left indent ---
$p1 = $d.Paragraphs(1)
$p1.Range.Borders.DistanceFromTop = 5
$p1.Range.Borders.DistanceFromLeft = 5
$p1.Range.Borders.DistanceFromBottom = 5
$p1.Range.Borders.DistanceFromRight = 5
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# --- Update the placeholder id text and drop the now-unneeded trailing space run ---
$findRange = $d.Content
[void]$findRange.Find.Execute("**ID__AFFARS_5311_topic_7__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$idStart = $findRange.Start
$idEnd = $findRange.End

# the single-space run immediately follows the placeholder text
$spaceRange = $d.Range($idEnd, $idEnd + 1)
if ($spaceRange.Text -eq " ") {
    $spaceRange.Delete()
}

$idRange = $d.Range($idStart, $idEnd)
$idRange.Text = "**ID__AFFARS_SUBPART_5311_5__ID**"
